$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'322.11"
$ws.Range("E2").Value = "'8.09%"
$ws.Range("G2").Value = "'7"

$ws.Range("D3").Value = "'49.13"
$ws.Range("E3").Value = "'17.67%"
$ws.Range("G3").Value = "'7"

$ws.Range("D4").Value = "'5.290"
$ws.Range("E4").Value = "'5.44%"
$ws.Range("G4").Value = "'7"

$ws.Range("D5").Value = "'0.08097"
$ws.Range("E5").Value = "'7.45%"
$ws.Range("G5").Value = "'7"

$ws.Range("D6").Value = "'4.617"
$ws.Range("E6").Value = "'5.48%"
$ws.Range("G6").Value = "'7"

$ws.Range("D7").Value = "'1.666"
$ws.Range("E7").Value = "'2.54%"
$ws.Range("G7").Value = "'7"

$ws.Range("D8").Value = "'1.219"
$ws.Range("E8").Value = "'32.27%"
$ws.Range("G8").Value = "'7"

$ws.Range("D9").Value = "'0.1334"
$ws.Range("E9").Value = "'12.92%"
$ws.Range("G9").Value = "'7"

$ws.Range("D10").Value = "'0.1963"
$ws.Range("E10").Value = "'7.46%"
$ws.Range("G10").Value = "'7"

$ws.Range("D11").Value = "'0.09528"
$ws.Range("E11").Value = "'6.44%"
$ws.Range("G11").Value = "'7"

$ws.Range("D12").Value = "'0.04509"
$ws.Range("E12").Value = "'10.52%"
$ws.Range("G12").Value = "'7"

$ws.Range("D13").Value = "'0.1047"
$ws.Range("E13").Value = "'-0.30%"
$ws.Range("G13").Value = "'7"

$ws.Range("D14").Value = "'0.001328"
$ws.Range("E14").Value = "'2.65%"
$ws.Range("G14").Value = "'7"

$ws.Range("D15").Value = "'0.005883"
$ws.Range("E15").Value = "'1.38%"
$ws.Range("G15").Value = "'7"

$ws.Range("D16").Value = "'3.363"
$ws.Range("E16").Value = "'0.69%"
$ws.Range("G16").Value = "'7"

$ws.Range("D17").Value = "'2.437"
$ws.Range("E17").Value = "'1.51%"
$ws.Range("G17").Value = "'7"

$ws.Range("D18").Value = "'0.3392"
$ws.Range("E18").Value = "'1.91%"
$ws.Range("G18").Value = "'7"

$ws.Range("D19").Value = "'8.218"
$ws.Range("E19").Value = "'-0.67%"
$ws.Range("G19").Value = "'7"

$ws.Range("D20").Value = "'0.1413"
$ws.Range("E20").Value = "'3.00%"
$ws.Range("G20").Value = "'7"

$ws.Range("D21").Value = "'0.2919"
$ws.Range("E21").Value = "'-9.34%"
$ws.Range("G21").Value = "'7"

$ws.Range("D22").Value = "'0.04307"
$ws.Range("E22").Value = "'5.22%"
$ws.Range("G22").Value = "'7"

$ws.Range("D23").Value = "'0.001310"
$ws.Range("E23").Value = "'3.38%"
$ws.Range("G23").Value = "'7"

$ws.Range("D24").Value = "'0.004251"
$ws.Range("E24").Value = "'9.20%"
$ws.Range("G24").Value = "'7"

$ws.Range("D25").Value = "'0.0001352"
$ws.Range("E25").Value = "'9.87%"
$ws.Range("G25").Value = "'7"

$ws.Range("D26").Value = "'0.0003548"
$ws.Range("E26").Value = "'-4.72%"
$ws.Range("G26").Value = "'7"

$ws.Range("G27").Value = "'7"

$ws.Range("G28").Value = "'7"

$ws.Range("G29").Value = "'7"

$ws.Range("G30").Value = "'7"

$ws.Range("G31").Value = "'7"

$ws.Range("G32").Value = "'7"

$ws.Range("G33").Value = "'7"

$ws.Range("G34").Value = "'7"

$ws.Range("G35").Value = "'7"

$ws.Range("G36").Value = "'7"

$ws.Range("G37").Value = "'7"

$ws.Range("D38").Value = "'0.02733"
$ws.Range("E38").Value = "'13.63%"
$ws.Range("G38").Value = "'7"

$ws.Range("D39").Value = "'0.05609"
$ws.Range("E39").Value = "'7.54%"
$ws.Range("G39").Value = "'7"

$ws.Range("D40").Value = "'0.006313"
$ws.Range("E40").Value = "'0.13%"
$ws.Range("G40").Value = "'7"

$ws.Range("D41").Value = "'0.007699"
$ws.Range("E41").Value = "'-1.60%"
$ws.Range("G41").Value = "'7"

$ws.Range("D42").Value = "'0.1444"
$ws.Range("E42").Value = "'8.94%"
$ws.Range("G42").Value = "'7"

$ws.Range("D43").Value = "'0.007706"
$ws.Range("E43").Value = "'4.24%"
$ws.Range("G43").Value = "'7"

$ws.Range("D44").Value = "'0.008101"
$ws.Range("E44").Value = "'6.34%"
$ws.Range("G44").Value = "'7"

$ws.Range("D45").Value = "'0.3194"
$ws.Range("E45").Value = "'7.62%"
$ws.Range("G45").Value = "'7"

$ws.Range("D46").Value = "'0.00006986"
$ws.Range("E46").Value = "'5.96%"
$ws.Range("G46").Value = "'7"

$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.13%"
$ws.Range("G47").Value = "'7"

$ws.Range("E48").Value = "'29.23%"
$ws.Range("G48").Value = "'7"

$ws.Range("D49").Value = "'0.004010"
$ws.Range("E49").Value = "'-4.59%"
$ws.Range("G49").Value = "'7"

$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.13%"
$ws.Range("G50").Value = "'7"

$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.13%"
$ws.Range("G51").Value = "'7"
